# Risk Management Log v2.3 - iteration 4 docs and meeting risks update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 20: clear the Risk Expiry Date (L20) that was mistakenly filled in ---
$ws.Range("L20").ClearContents()

# --- Row 23: fill in the new risk entry (email bounce / validation risk) ---
$ws.Range("B23").Value = 20
$ws.Range("C23").Formula = '=IF(H23<=7, "Low Risk",IF(H23>=17,"High Risk","Medium Risk"))'
$ws.Range("D23").Value = "Email only allows 10 bounces a day "
$ws.Range("E23").Value = "Charlotte Hutchinson"
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 4
$ws.Range("H23").Formula = "=F23*G23"
$ws.Range("I23").Value = "External Influences"
$ws.Range("J23").Value = "Retain "
$ws.Range("K23").Value = "We validate all emails using javaScript"
$ws.Range("L23").Value = 41992

# --- Update the sheet view: scroll position and active selection ---
$ws.Activate()
$ws.Range("K24").Select()
